$wb = $excel.ActiveWorkbook

# --- Sheet: Vehicle Sections ---
$ws1 = $wb.Worksheets.Item("Vehicle Sections")

# Length (ft) column updates for Nosecone, Recovery Bay, Helium Bay, Upper Airframe
$ws1.Range("C2").Value = 0.5
$ws1.Range("C3").Value = 0.5
$ws1.Range("C4").Value = 0.5
$ws1.Range("C5").Value = 0.5

# Lox Tank row (row 6): Mass (lbs) and Length (ft, now a formula)
$ws1.Range("B6").Value = 4.4242043000000004
$ws1.Range("C6").Formula = "=4.5737205 / 12"

# Mid Airframe row (row 7): Length (ft)
$ws1.Range("C7").Value = 0.23774000000000001

# Fuel Tank row (row 8): Mass (lbs) and Length (ft, now a formula)
$ws1.Range("B8").Value = 4.4242043000000004
$ws1.Range("C8").Formula = "=4.5737205 / 12"

# Lower Airframe row (row 9): Length (ft)
$ws1.Range("C9").Value = 2.4

# Engine row (row 10): Length (ft)
$ws1.Range("C10").Value = 2.4

# Update selection on this sheet to C10
$ws1.Activate() | Out-Null
$ws1.Range("C10").Select() | Out-Null

# --- Sheet: Aerodynamic Properties ---
$ws2 = $wb.Worksheets.Item("Aerodynamic Properties")
$ws2.Activate() | Out-Null
$ws2.Range("E33").Select() | Out-Null

# Re-activate Vehicle Sections as the last active/selected sheet (matches tabSelected in diff)
$ws1.Activate() | Out-Null

$wb.Save()
